$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.060.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = "'3.336.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = "'583.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = "'175.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.90%  '
$ws.Range("D9").Value = "'0.182"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.06%  '
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("D11").Value = "'47.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.72%  '
$ws.Range("E12").Value = '  +1.55%  '
$ws.Range("D13").Value = "'701.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.55%  '
$ws.Range("D14").Value = "'3.875.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").Value = "'8.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = "'68.054.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = "'3.317.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").Value = "'17.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("E20").Value = '  +2.53%  '
$ws.Range("E21").Value = '  +0.85%  '
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = "'16.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").Value = "'100.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.40%  '
$ws.Range("D25").Value = "'3.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("D26").Value = "'2.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = "'9.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.10%  '
$ws.Range("D28").Value = "'33.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").Value = "'8.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("D30").Value = "'7.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.73%  '
$ws.Range("D31").Value = "'569.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("D32").Value = "'11.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  +1.57%  '
$ws.Range("D34").Value = "'3.755.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = "'57.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.99%  '
$ws.Range("D37").Value = "'3.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").Value = "'35.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.69%  '
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").Value = "'0.0₃0677"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.37%  '
$ws.Range("D43").Value = "'0.334"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("D44").Value = "'3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("D46").Value = "'2.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.02%  '
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("D50").Value = "'130.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("E51").Value = '  +0.27%  '
